$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 98, shifting the existing rows 98:104 down to 99:105
$ws.Rows.Item(98).Insert()

# Populate the newly inserted row 98 with the new weekly price record
$ws.Range("A98").Value = 11
$ws.Range("B98").Value = "Vega Monumental Concepción"
$ws.Range("C98").Value = "Bíobío"
$ws.Range("D98").Value = 45021
$ws.Range("E98").Value = 8
$ws.Range("F98").Value = 100112037
$ws.Range("G98").Value = "Cebollín"
$ws.Range("H98").Value = "Sin especificar"
$ws.Range("I98").Value = "Primera"
$ws.Range("J98").Value = 60
$ws.Range("K98").Value = 5000
$ws.Range("L98").Value = 5500
$ws.Range("M98").Value = 5250
$ws.Range("N98").Value = "$/paquete 36 unidades"
$ws.Range("O98").Value = "Región Metropolitana"
$ws.Range("P98").Value = 146
$ws.Range("Q98").Value = 36
$ws.Range("R98").Value = "Hortaliza"
